$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 0.1845238095238095
$ws.Cells.Item(2, 3).Value = 0.5982142857142857
$ws.Cells.Item(2, 10).Value = 0.02380952380952381
$ws.Cells.Item(2, 16).Value = 0.1220238095238095
$ws.Cells.Item(2, 19).Value = 0.07142857142857142
$ws.Cells.Item(3, 3).Value = 0.005
$ws.Cells.Item(3, 10).Value = 0.075
$ws.Cells.Item(3, 16).Value = 0.73
$ws.Cells.Item(3, 19).Value = 0.19
$ws.Cells.Item(4, 10).Value = 0.1142857142857143
$ws.Cells.Item(4, 16).Value = 0.6857142857142857
$ws.Cells.Item(4, 19).Value = 0.2
$ws.Cells.Item(6, 2).Value = 0.04878048780487805
$ws.Cells.Item(6, 4).Value = 0.008130081300813009
$ws.Cells.Item(6, 6).Value = 0.06504065040650407
$ws.Cells.Item(6, 10).Value = 0.2682926829268293
$ws.Cells.Item(6, 15).Value = 0.02032520325203252
$ws.Cells.Item(6, 17).Value = 0.1747967479674797
$ws.Cells.Item(6, 18).Value = 0.08536585365853659
$ws.Cells.Item(6, 19).Value = 0.3292682926829268
$ws.Cells.Item(7, 2).Value = 0.09473684210526316
$ws.Cells.Item(7, 4).Value = 0.005263157894736842
$ws.Cells.Item(7, 5).Value = 0.005263157894736842
$ws.Cells.Item(7, 6).Value = 0.04736842105263158
$ws.Cells.Item(7, 10).Value = 0.131578947368421
$ws.Cells.Item(7, 15).Value = 0.02105263157894737
$ws.Cells.Item(7, 17).Value = 0.1684210526315789
$ws.Cells.Item(7, 18).Value = 0.1210526315789474
$ws.Cells.Item(7, 19).Value = 0.4052631578947368
$ws.Cells.Item(8, 2).Value = 0.09513742071881606
$ws.Cells.Item(8, 4).Value = 0.01691331923890063
$ws.Cells.Item(8, 5).Value = 0.002114164904862579
$ws.Cells.Item(8, 6).Value = 0.06553911205073996
$ws.Cells.Item(8, 10).Value = 0.1099365750528541
$ws.Cells.Item(8, 15).Value = 0.02114164904862579
$ws.Cells.Item(8, 17).Value = 0.1627906976744186
$ws.Cells.Item(8, 18).Value = 0.105708245243129
$ws.Cells.Item(8, 19).Value = 0.4207188160676533
$ws.Cells.Item(9, 2).Value = 0.07111111111111111
$ws.Cells.Item(9, 4).Value = 0.008888888888888889
$ws.Cells.Item(9, 6).Value = 0.05333333333333334
$ws.Cells.Item(9, 10).Value = 0.1022222222222222
$ws.Cells.Item(9, 15).Value = 0.008888888888888889
$ws.Cells.Item(9, 17).Value = 0.1688888888888889
$ws.Cells.Item(9, 18).Value = 0.1377777777777778
$ws.Cells.Item(9, 19).Value = 0.4488888888888889
$ws.Cells.Item(10, 2).Value = 0.1200269723533378
$ws.Cells.Item(10, 4).Value = 0.01618341200269724
$ws.Cells.Item(10, 5).Value = 0.0006743088334457181
$ws.Cells.Item(10, 6).Value = 0.06540795684423466
$ws.Cells.Item(10, 10).Value = 0.1362103843560351
$ws.Cells.Item(10, 15).Value = 0.01753202966958867
$ws.Cells.Item(10, 17).Value = 0.2016183412002697
$ws.Cells.Item(10, 18).Value = 0.09912339851652056
$ws.Cells.Item(10, 19).Value = 0.3432231962238705
$ws.Cells.Item(11, 7).Value = 0.1493055555555556
$ws.Cells.Item(11, 10).Value = 0.06944444444444445
$ws.Cells.Item(11, 11).Value = 0.2013888888888889
$ws.Cells.Item(11, 12).Value = 0.5659722222222222
$ws.Cells.Item(11, 19).Value = 0.01388888888888889
$ws.Cells.Item(12, 7).Value = 0.7062146892655368
$ws.Cells.Item(12, 10).Value = 0.2203389830508475
$ws.Cells.Item(12, 11).Value = 0.01694915254237288
$ws.Cells.Item(12, 12).Value = 0.03954802259887006
$ws.Cells.Item(12, 19).Value = 0.01694915254237288
$ws.Cells.Item(13, 7).Value = 0.7380952380952381
$ws.Cells.Item(13, 10).Value = 0.2380952380952381
$ws.Cells.Item(13, 19).Value = 0.02380952380952381
$ws.Cells.Item(14, 7).Value = 0.6
$ws.Cells.Item(14, 10).Value = 0.4
$ws.Cells.Item(15, 8).Value = 0.141025641025641
$ws.Cells.Item(15, 9).Value = 0.05555555555555555
$ws.Cells.Item(15, 10).Value = 0.4145299145299146
$ws.Cells.Item(15, 11).Value = 0.04700854700854701
$ws.Cells.Item(15, 13).Value = 0.008547008547008548
$ws.Cells.Item(15, 14).Value = 0.004273504273504274
$ws.Cells.Item(15, 15).Value = 0.0811965811965812
$ws.Cells.Item(15, 19).Value = 0.2478632478632479
$ws.Cells.Item(16, 6).Value = 0.01951219512195122
$ws.Cells.Item(16, 8).Value = 0.1170731707317073
$ws.Cells.Item(16, 9).Value = 0.07804878048780488
$ws.Cells.Item(16, 10).Value = 0.5073170731707317
$ws.Cells.Item(16, 11).Value = 0.05853658536585366
$ws.Cells.Item(16, 13).Value = 0.02439024390243903
$ws.Cells.Item(16, 14).Value = 0.004878048780487805
$ws.Cells.Item(16, 15).Value = 0.05853658536585366
$ws.Cells.Item(16, 19).Value = 0.1317073170731707
$ws.Cells.Item(17, 6).Value = 0.02904564315352697
$ws.Cells.Item(17, 8).Value = 0.2074688796680498
$ws.Cells.Item(17, 9).Value = 0.09336099585062241
$ws.Cells.Item(17, 10).Value = 0.3672199170124482
$ws.Cells.Item(17, 11).Value = 0.08921161825726141
$ws.Cells.Item(17, 13).Value = 0.02074688796680498
$ws.Cells.Item(17, 14).Value = 0.002074688796680498
$ws.Cells.Item(17, 15).Value = 0.06639004149377593
$ws.Cells.Item(17, 19).Value = 0.1244813278008299
$ws.Cells.Item(18, 6).Value = 0.02592592592592593
$ws.Cells.Item(18, 8).Value = 0.1666666666666667
$ws.Cells.Item(18, 9).Value = 0.08888888888888889
$ws.Cells.Item(18, 10).Value = 0.4703703703703704
$ws.Cells.Item(18, 11).Value = 0.1
$ws.Cells.Item(18, 13).Value = 0.007407407407407408
$ws.Cells.Item(18, 14).Value = 0.007407407407407408
$ws.Cells.Item(18, 15).Value = 0.02962962962962963
$ws.Cells.Item(18, 19).Value = 0.1037037037037037
$ws.Cells.Item(19, 6).Value = 0.01597676107480029
$ws.Cells.Item(19, 8).Value = 0.1902687000726216
$ws.Cells.Item(19, 9).Value = 0.0915032679738562
$ws.Cells.Item(19, 10).Value = 0.3907044299201162
$ws.Cells.Item(19, 11).Value = 0.09368191721132897
$ws.Cells.Item(19, 13).Value = 0.01742919389978214
$ws.Cells.Item(19, 15).Value = 0.06681190994916485
$ws.Cells.Item(19, 19).Value = 0.1336238198983297
